$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.048.63'
$ws.Range("E2").Value = '  +6.20%  '

$ws.Range("D3").Value = '3.111.85'
$ws.Range("E3").Value = '  +3.81%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '''587.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.59%  '

$ws.Range("D6").Value = '''144.01'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.59%  '

$ws.Range("E7").Value = '  -0.10%  '

$ws.Range("D8").Value = '3.100.77'
$ws.Range("E8").Value = '  +3.88%  '

$ws.Range("E9").Value = '  +2.50%  '

$ws.Range("E10").Value = '  +10.11%  '

$ws.Range("D11").Value = '''5.69'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +9.34%  '

$ws.Range("E12").Value = '  +1.97%  '

$ws.Range("D13").Value = '''0.0000245'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.97%  '

$ws.Range("D14").Value = '''35.59'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.35%  '

$ws.Range("E15").Value = '  +0.81%  '

$ws.Range("D16").Value = '3.627.97'
$ws.Range("E16").Value = '  +3.87%  '

$ws.Range("D17").Value = '''7.27'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.41%  '

$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '62.983.49'
$ws.Range("E18").Value = '  +6.12%  '

$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value = '3.111.68'
$ws.Range("E19").Value = '  +4.13%  '

$ws.Range("D20").Value = '''454.14'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.41%  '

$ws.Range("E21").Value = '  +3.33%  '

$ws.Range("D22").Value = '''0.734'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.68%  '

$ws.Range("D23").Value = '''7.61'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.54%  '

$ws.Range("D24").Value = '''13.67'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.99%  '

$ws.Range("D25").Value = '''82.10'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.30%  '

$ws.Range("E26").Value = '  +0.12%  '

$ws.Range("E27").Value = '  +0.94%  '

$ws.Range("D28").Value = '''2.72'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.31%  '

$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").Value = '''8.31'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.16%  '

$ws.Range("B30").Value = 'FirstDigitalUSD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D30").Value = '''1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.04%  '

$ws.Range("D31").Value = '''6.90'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +12.47%  '

$ws.Range("E32").Value = '  +11.09%  '

$ws.Range("D33").Value = '''27.11'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.23%  '

$ws.Range("B34").Value = 'PEPE'
$ws.Range("C34").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D34").Value = '0.0₃0817'
$ws.Range("E34").Value = '  +7.20%  '

$ws.Range("B35").Value = 'Stacks'
$ws.Range("C35").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D35").Value = '''2.33'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +9.95%  '

$ws.Range("B36").Value = 'Mantle'
$ws.Range("C36").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D36").Value = '''1.04'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.25%  '

$ws.Range("D37").Value = '''6.07'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.60%  '

$ws.Range("B38").Value = 'dogwifhat'
$ws.Range("C38").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D38").Value = '''3.07'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +11.36%  '

$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").Value = '''50.97'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.16%  '

$ws.Range("E40").Value = '  +1.17%  '

$ws.Range("D41").Value = '''430.15'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.27%  '

$ws.Range("D42").Value = '2.965.99'
$ws.Range("E42").Value = '  +6.70%  '

$ws.Range("E43").Value = '  +6.04%  '

$ws.Range("E44").Value = '  +3.62%  '

$ws.Range("D45").Value = '''0.275'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +8.79%  '

$ws.Range("E46").Value = '  +7.56%  '

$ws.Range("D47").Value = '''124.99'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.10%  '

$ws.Range("D49").Value = '''34.65'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.69%  '

$ws.Range("E50").Value = '  +1.01%  '

$ws.Range("D51").Value = '''24.85'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.39%  '
